# Apply KHL stats update: append 2025-11-16 games and recompute shots-on-goal aggregates.
# (chore(runtime): publish files + archive (2025-11-17 11:04:16))

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Matches_SOG: append the 7 games played 2025-11-16 (rows 477-483)
# ---------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$wsMatches.Cells.Item(477, 1).NumberFormat = "@"
$wsMatches.Cells.Item(477, 1).Value = "897775"
$wsMatches.Cells.Item(477, 2).Value = "2025-11-16T10:00:00"
$wsMatches.Cells.Item(477, 3).Value = "Амур"
$wsMatches.Cells.Item(477, 4).Value = "Адмирал"
$wsMatches.Cells.Item(477, 5).Value = 23
$wsMatches.Cells.Item(477, 6).Value = 26
$wsMatches.Cells.Item(477, 7).Value = "khl_text"

$wsMatches.Cells.Item(478, 1).NumberFormat = "@"
$wsMatches.Cells.Item(478, 1).Value = "897774"
$wsMatches.Cells.Item(478, 2).Value = "2025-11-16T14:00:00"
$wsMatches.Cells.Item(478, 3).Value = "Трактор"
$wsMatches.Cells.Item(478, 4).Value = "Автомобилист"
$wsMatches.Cells.Item(478, 5).Value = 40
$wsMatches.Cells.Item(478, 6).Value = 29
$wsMatches.Cells.Item(478, 7).Value = "khl_text"

$wsMatches.Cells.Item(479, 1).NumberFormat = "@"
$wsMatches.Cells.Item(479, 1).Value = "897776"
$wsMatches.Cells.Item(479, 2).Value = "2025-11-16T16:00:00"
$wsMatches.Cells.Item(479, 3).Value = "Лада"
$wsMatches.Cells.Item(479, 4).Value = "Нефтехимик"
$wsMatches.Cells.Item(479, 5).Value = 31
$wsMatches.Cells.Item(479, 6).Value = 41
$wsMatches.Cells.Item(479, 7).Value = "khl_text"

$wsMatches.Cells.Item(480, 1).NumberFormat = "@"
$wsMatches.Cells.Item(480, 1).Value = "897777"
$wsMatches.Cells.Item(480, 2).Value = "2025-11-16T17:00:00"
$wsMatches.Cells.Item(480, 3).Value = "Северсталь"
$wsMatches.Cells.Item(480, 4).Value = "ЦСКА"
$wsMatches.Cells.Item(480, 5).Value = 22
$wsMatches.Cells.Item(480, 6).Value = 27
$wsMatches.Cells.Item(480, 7).Value = "khl_text"

$wsMatches.Cells.Item(481, 1).NumberFormat = "@"
$wsMatches.Cells.Item(481, 1).Value = "897778"
$wsMatches.Cells.Item(481, 2).Value = "2025-11-16T17:00:00"
$wsMatches.Cells.Item(481, 3).Value = "Ак Барс"
$wsMatches.Cells.Item(481, 4).Value = "Динамо Мн"
$wsMatches.Cells.Item(481, 5).Value = 34
$wsMatches.Cells.Item(481, 6).Value = 27
$wsMatches.Cells.Item(481, 7).Value = "khl_text"

$wsMatches.Cells.Item(482, 1).NumberFormat = "@"
$wsMatches.Cells.Item(482, 1).Value = "897779"
$wsMatches.Cells.Item(482, 2).Value = "2025-11-16T17:00:00"
$wsMatches.Cells.Item(482, 3).Value = "Торпедо"
$wsMatches.Cells.Item(482, 4).Value = "Авангард"
$wsMatches.Cells.Item(482, 5).Value = 43
$wsMatches.Cells.Item(482, 6).Value = 43
$wsMatches.Cells.Item(482, 7).Value = "khl_text"

$wsMatches.Cells.Item(483, 1).NumberFormat = "@"
$wsMatches.Cells.Item(483, 1).Value = "897780"
$wsMatches.Cells.Item(483, 2).Value = "2025-11-16T17:30:00"
$wsMatches.Cells.Item(483, 3).Value = "СКА"
$wsMatches.Cells.Item(483, 4).Value = "Металлург Мг"
$wsMatches.Cells.Item(483, 5).Value = 22
$wsMatches.Cells.Item(483, 6).Value = 35
$wsMatches.Cells.Item(483, 7).Value = "khl_text"

# ---------------------------------------------------------------
# 2) Shots_HA: refresh as_of_utc for every team + HOGF/HOGA/AOGF/AOGA
#    totals & per-game rates for the 14 teams that played on 2025-11-16
# ---------------------------------------------------------------
$wsShotsHA = $wb.Worksheets.Item("Shots_HA")

$wsShotsHA.Cells.Item(2, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(2, 6).Value = 19
$wsShotsHA.Cells.Item(2, 11).Value = 685
$wsShotsHA.Cells.Item(2, 12).Value = 593
$wsShotsHA.Cells.Item(2, 13).Value = 36.1
$wsShotsHA.Cells.Item(2, 14).Value = 31.2

$wsShotsHA.Cells.Item(3, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(3, 6).Value = 28
$wsShotsHA.Cells.Item(3, 11).Value = 781
$wsShotsHA.Cells.Item(3, 12).Value = 871
$wsShotsHA.Cells.Item(3, 14).Value = 31.1

$wsShotsHA.Cells.Item(4, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(4, 6).Value = 21
$wsShotsHA.Cells.Item(4, 11).Value = 664
$wsShotsHA.Cells.Item(4, 12).Value = 584
$wsShotsHA.Cells.Item(4, 13).Value = 31.6
$wsShotsHA.Cells.Item(4, 14).Value = 27.8

$wsShotsHA.Cells.Item(5, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(5, 5).Value = 25
$wsShotsHA.Cells.Item(5, 7).Value = 846
$wsShotsHA.Cells.Item(5, 8).Value = 641

$wsShotsHA.Cells.Item(6, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(6, 5).Value = 22
$wsShotsHA.Cells.Item(6, 7).Value = 664
$wsShotsHA.Cells.Item(6, 8).Value = 772
$wsShotsHA.Cells.Item(6, 9).Value = 30.2
$wsShotsHA.Cells.Item(6, 10).Value = 35.1

$wsShotsHA.Cells.Item(7, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsHA.Cells.Item(8, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsHA.Cells.Item(9, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(9, 6).Value = 19
$wsShotsHA.Cells.Item(9, 11).Value = 695
$wsShotsHA.Cells.Item(9, 12).Value = 518
$wsShotsHA.Cells.Item(9, 13).Value = 36.6
$wsShotsHA.Cells.Item(9, 14).Value = 27.3

$wsShotsHA.Cells.Item(10, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsHA.Cells.Item(11, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(11, 5).Value = 23
$wsShotsHA.Cells.Item(11, 7).Value = 629
$wsShotsHA.Cells.Item(11, 8).Value = 817
$wsShotsHA.Cells.Item(11, 9).Value = 27.3
$wsShotsHA.Cells.Item(11, 10).Value = 35.5

$wsShotsHA.Cells.Item(12, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsHA.Cells.Item(13, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(13, 6).Value = 18
$wsShotsHA.Cells.Item(13, 11).Value = 523
$wsShotsHA.Cells.Item(13, 12).Value = 482
$wsShotsHA.Cells.Item(13, 13).Value = 29.1
$wsShotsHA.Cells.Item(13, 14).Value = 26.8

$wsShotsHA.Cells.Item(14, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(14, 6).Value = 20
$wsShotsHA.Cells.Item(14, 11).Value = 558
$wsShotsHA.Cells.Item(14, 12).Value = 753
$wsShotsHA.Cells.Item(14, 13).Value = 27.9
$wsShotsHA.Cells.Item(14, 14).Value = 37.6

$wsShotsHA.Cells.Item(15, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(15, 5).Value = 26
$wsShotsHA.Cells.Item(15, 7).Value = 852
$wsShotsHA.Cells.Item(15, 8).Value = 869
$wsShotsHA.Cells.Item(15, 9).Value = 32.8

$wsShotsHA.Cells.Item(16, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsHA.Cells.Item(17, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(17, 5).Value = 17
$wsShotsHA.Cells.Item(17, 7).Value = 491
$wsShotsHA.Cells.Item(17, 8).Value = 390
$wsShotsHA.Cells.Item(17, 9).Value = 28.9
$wsShotsHA.Cells.Item(17, 10).Value = 22.9

$wsShotsHA.Cells.Item(18, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsHA.Cells.Item(19, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsHA.Cells.Item(20, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(20, 5).Value = 23
$wsShotsHA.Cells.Item(20, 7).Value = 760
$wsShotsHA.Cells.Item(20, 8).Value = 693
$wsShotsHA.Cells.Item(20, 9).Value = 33.0
$wsShotsHA.Cells.Item(20, 10).Value = 30.1

$wsShotsHA.Cells.Item(21, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(21, 5).Value = 19
$wsShotsHA.Cells.Item(21, 7).Value = 639
$wsShotsHA.Cells.Item(21, 8).Value = 570
$wsShotsHA.Cells.Item(21, 9).Value = 33.6
$wsShotsHA.Cells.Item(21, 10).Value = 30.0

$wsShotsHA.Cells.Item(22, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsHA.Cells.Item(23, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsHA.Cells.Item(23, 6).Value = 23
$wsShotsHA.Cells.Item(23, 11).Value = 578
$wsShotsHA.Cells.Item(23, 12).Value = 646
$wsShotsHA.Cells.Item(23, 13).Value = 25.1
$wsShotsHA.Cells.Item(23, 14).Value = 28.1

# ---------------------------------------------------------------
# 3) Shots_Summary: refresh as_of_utc + GP/SOG/SOGA totals & per-game
#    rates (derived from the Shots_HA home+away totals above)
# ---------------------------------------------------------------
$wsShotsSummary = $wb.Worksheets.Item("Shots_Summary")

$wsShotsSummary.Cells.Item(2, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(2, 5).Value = 42
$wsShotsSummary.Cells.Item(2, 6).Value = 1436
$wsShotsSummary.Cells.Item(2, 7).Value = 1244
$wsShotsSummary.Cells.Item(2, 8).Value = 34.2
$wsShotsSummary.Cells.Item(2, 9).Value = 29.6

$wsShotsSummary.Cells.Item(3, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(3, 5).Value = 46
$wsShotsSummary.Cells.Item(3, 6).Value = 1302
$wsShotsSummary.Cells.Item(3, 7).Value = 1429
$wsShotsSummary.Cells.Item(3, 9).Value = 31.1

$wsShotsSummary.Cells.Item(4, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(4, 5).Value = 39
$wsShotsSummary.Cells.Item(4, 6).Value = 1350
$wsShotsSummary.Cells.Item(4, 7).Value = 1067
$wsShotsSummary.Cells.Item(4, 8).Value = 34.6
$wsShotsSummary.Cells.Item(4, 9).Value = 27.4

$wsShotsSummary.Cells.Item(5, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(5, 5).Value = 46
$wsShotsSummary.Cells.Item(5, 6).Value = 1557
$wsShotsSummary.Cells.Item(5, 7).Value = 1268

$wsShotsSummary.Cells.Item(6, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(6, 5).Value = 43
$wsShotsSummary.Cells.Item(6, 6).Value = 1256
$wsShotsSummary.Cells.Item(6, 7).Value = 1541
$wsShotsSummary.Cells.Item(6, 8).Value = 29.2
$wsShotsSummary.Cells.Item(6, 9).Value = 35.8

$wsShotsSummary.Cells.Item(7, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsSummary.Cells.Item(8, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsSummary.Cells.Item(9, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(9, 5).Value = 44
$wsShotsSummary.Cells.Item(9, 6).Value = 1615
$wsShotsSummary.Cells.Item(9, 7).Value = 1195
$wsShotsSummary.Cells.Item(9, 8).Value = 36.7
$wsShotsSummary.Cells.Item(9, 9).Value = 27.2

$wsShotsSummary.Cells.Item(10, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsSummary.Cells.Item(11, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(11, 5).Value = 44
$wsShotsSummary.Cells.Item(11, 6).Value = 1160
$wsShotsSummary.Cells.Item(11, 7).Value = 1614
$wsShotsSummary.Cells.Item(11, 8).Value = 26.4
$wsShotsSummary.Cells.Item(11, 9).Value = 36.7

$wsShotsSummary.Cells.Item(12, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsSummary.Cells.Item(13, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(13, 5).Value = 44
$wsShotsSummary.Cells.Item(13, 6).Value = 1449
$wsShotsSummary.Cells.Item(13, 7).Value = 1142

$wsShotsSummary.Cells.Item(14, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(14, 5).Value = 46
$wsShotsSummary.Cells.Item(14, 6).Value = 1366
$wsShotsSummary.Cells.Item(14, 7).Value = 1639
$wsShotsSummary.Cells.Item(14, 8).Value = 29.7
$wsShotsSummary.Cells.Item(14, 9).Value = 35.6

$wsShotsSummary.Cells.Item(15, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(15, 5).Value = 44
$wsShotsSummary.Cells.Item(15, 6).Value = 1406
$wsShotsSummary.Cells.Item(15, 7).Value = 1459
$wsShotsSummary.Cells.Item(15, 8).Value = 32.0
$wsShotsSummary.Cells.Item(15, 9).Value = 33.2

$wsShotsSummary.Cells.Item(16, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsSummary.Cells.Item(17, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(17, 5).Value = 43
$wsShotsSummary.Cells.Item(17, 6).Value = 1345
$wsShotsSummary.Cells.Item(17, 7).Value = 1082
$wsShotsSummary.Cells.Item(17, 8).Value = 31.3
$wsShotsSummary.Cells.Item(17, 9).Value = 25.2

$wsShotsSummary.Cells.Item(18, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsSummary.Cells.Item(19, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsSummary.Cells.Item(20, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(20, 5).Value = 50
$wsShotsSummary.Cells.Item(20, 6).Value = 1691
$wsShotsSummary.Cells.Item(20, 7).Value = 1559
$wsShotsSummary.Cells.Item(20, 8).Value = 33.8
$wsShotsSummary.Cells.Item(20, 9).Value = 31.2

$wsShotsSummary.Cells.Item(21, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(21, 5).Value = 46
$wsShotsSummary.Cells.Item(21, 6).Value = 1562
$wsShotsSummary.Cells.Item(21, 7).Value = 1438
$wsShotsSummary.Cells.Item(21, 8).Value = 34.0

$wsShotsSummary.Cells.Item(22, 4).Value = "2025-11-16T17:30:00Z"

$wsShotsSummary.Cells.Item(23, 4).Value = "2025-11-16T17:30:00Z"
$wsShotsSummary.Cells.Item(23, 5).Value = 43
$wsShotsSummary.Cells.Item(23, 6).Value = 1048
$wsShotsSummary.Cells.Item(23, 7).Value = 1223
$wsShotsSummary.Cells.Item(23, 8).Value = 24.4
$wsShotsSummary.Cells.Item(23, 9).Value = 28.4

# ---------------------------------------------------------------
# 4) Meta_ext: bump as_of_utc + build_version
# ---------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Cells.Item(2, 2).Value = "2025-11-16T17:30:00Z"
$wsMeta.Cells.Item(2, 4).Value = 68
